$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update scoring formulas/labels for CAPE, AQ and BISBAS rows
$ws.Range("D3").Value = "mean(distress * frequency)"
$ws.Range("D4").Value = "sum(agree)"

# BISBAS row: min changed from 0 to 1, mode label changed
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = "sum(bas only)"

# Auto-fit column D to match the new, wider label text
$ws.Columns.Item(4).ColumnWidth = 24.59

# Move the active selection to D7, matching the last-edited cell
$ws.Range("D7").Select() | Out-Null

$wb.Save()
